$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.084.77"
$ws.Range("E2").Value = "  -2.71%  "

# Row 3
$ws.Range("D3").Value = "1.733.08"
$ws.Range("E3").Value = "  -1.36%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Value = "'310.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.18%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7
$ws.Range("D7").Value = "'0.4884"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.12%  "

# Row 8
$ws.Range("E8").Value = "  +0.64%  "

# Row 9
$ws.Range("D9").Value = "'43.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.37%  "

# Row 10
$ws.Range("D10").Value = "'0.07299"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "

# Row 11
$ws.Range("D11").Value = "'1.055"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.69%  "

# Row 12
$ws.Range("D12").Value = "'0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("E13").Value = "  -2.87%  "

# Row 14
$ws.Range("D14").Value = "'5.900"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.32%  "

# Row 15
$ws.Range("D15").Value = "1.733.94"
$ws.Range("E15").Value = "  -1.32%  "

# Row 16
$ws.Range("D16").Value = "'6.911"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.68%  "

# Row 17
$ws.Range("D17").Value = "'87.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.51%  "

# Row 18
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "

# Row 19
$ws.Range("D19").Value = "'0.06416"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("D20").Value = "'0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("D21").Value = "'16.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.30%  "

# Row 22
$ws.Range("D22").Value = "'5.720"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
$ws.Range("D23").Value = "27.127.44"
$ws.Range("E23").Value = "  -2.67%  "

# Row 24
$ws.Range("D24").Value = "'10.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "

# Row 25
$ws.Range("D25").Value = "'2.082"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.58%  "

# Row 26
$ws.Range("D26").Value = "'154.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.97%  "

# Row 27
$ws.Range("D27").Value = "'20.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("D28").Value = "1.924.75"
$ws.Range("E28").Value = "  -1.78%  "

# Row 29
$ws.Range("D29").Value = "'2.092"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.04%  "

# Row 30
$ws.Range("D30").Value = "'121.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.35%  "

# Row 31
$ws.Range("D31").Value = "'1.053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.20%  "

# Row 32
$ws.Range("D32").Value = "'0.09333"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "

# Row 33
$ws.Range("D33").Value = "'3.648"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34
$ws.Range("E34").Value = "  -1.90%  "

# Row 35
$ws.Range("D35").Value = "'0.05952"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.31%  "

# Row 36
$ws.Range("E36").Value = "  -2.50%  "

# Row 37
$ws.Range("E37").Value = "  -6.00%  "

# Row 38
$ws.Range("D38").Value = "'1.431"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.65%  "

# Row 39
$ws.Range("D39").Value = "'4.794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.16%  "

# Row 40
$ws.Range("D40").Value = "'0.2002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41
$ws.Range("D41").Value = "'0.6019"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.86%  "

# Row 42
$ws.Range("E42").Value = "  -0.07%  "

# Row 43
$ws.Range("D43").Value = "'1.099"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.74%  "

# Row 44
$ws.Range("D44").Value = "'7.498"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.32%  "

# Row 45
$ws.Range("D45").Value = "'12.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.08%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5689"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.09%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.587"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "

# Row 48
$ws.Range("D48").Value = "'118.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.82%  "

# Row 49
$ws.Range("D49").Value = "'1.854"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.76%  "

# Row 50
$ws.Range("D50").Value = "'1.108"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.13%  "

# Row 51
$ws.Range("D51").Value = "'0.06651"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
